# Handles float input without breaking stuff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Summary block (rows 10-12): scoring numbers updated, and row-label cells
# (A10/A11/A12) adopt the same "header" style already used by A9 (s=4).
# ---------------------------------------------------------------------------
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B10").Value = 11
$ws.Range("C10").Value = 6
$ws.Range("D10").Value = 11
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("B12").Value = 44
$ws.Range("C12").Value = -6
$ws.Range("E12").Value = "38/112"

# ---------------------------------------------------------------------------
# Third answer block (columns G/H) is gone entirely; most of the second
# answer block (columns D/E, rows 19-40) is gone too - only D16:E18 remain,
# with D16:D18 repurposed to hold the (now shorter) student-answer list.
# ---------------------------------------------------------------------------
$ws.Range("G15:H40").Clear()
$ws.Range("D19:E40").Clear()

# ---------------------------------------------------------------------------
# Student-answer cells: blank "Student Ans" cells (style 7, black) become
# filled in with the option the student actually picked. Matching answers
# get the "correct" green style (copied from B10, which already carries it);
# non-matching answers get the "incorrect" red style (copied from C11).
# ---------------------------------------------------------------------------
$ws.Range("B10").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("A27").PasteSpecial(-4122)
$ws.Range("A32").PasteSpecial(-4122)
$ws.Range("A33").PasteSpecial(-4122)
$ws.Range("A35").PasteSpecial(-4122)
$ws.Range("A36").PasteSpecial(-4122)
$ws.Range("A38").PasteSpecial(-4122)
$ws.Range("A39").PasteSpecial(-4122)
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C11").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("A25").PasteSpecial(-4122)
$ws.Range("A29").PasteSpecial(-4122)
$ws.Range("A37").PasteSpecial(-4122)
$ws.Range("D17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A16").Value = "Option A"
$ws.Range("A18").Value = "Option C"
$ws.Range("A19").Value = "Option B"
$ws.Range("A21").Value = "Option C"
$ws.Range("A25").Value = "Option B"
$ws.Range("A27").Value = "Option A"
$ws.Range("A29").Value = "Option C"
$ws.Range("A32").Value = "Option C"
$ws.Range("A33").Value = "Option D"
$ws.Range("A35").Value = "Option D"
$ws.Range("A36").Value = "Option A"
$ws.Range("A37").Value = "Option C"
$ws.Range("A38").Value = "Option A"
$ws.Range("A39").Value = "Option D"

$ws.Range("D16").Value = "Option A"
$ws.Range("D17").Value = "Option B"
$ws.Range("D18").Value = "Option D"
